# "Add levels up to 50" - populate additional level/enemy data rows.
$wb = $excel.ActiveWorkbook

$wsMeta    = $wb.Worksheets.Item("meta")
$ws2       = $wb.Worksheets.Item("level")
$ws3       = $wb.Worksheets.Item("enemies")
$ws4       = $wb.Worksheets.Item("misc")

# -----------------------------------------------------------------
# Sheet "misc": add the new row 5 entry (was previously a gap in the
# list). This introduces the new shared string "...Butterfly Arrow"
# which - chronologically - is the first brand-new string written,
# so it must happen before the new "enemies" sheet strings below.
# -----------------------------------------------------------------
$ws4.Range("A5").Value = "EnemyPrefabs/Arrow Enemies/Butterfly/Butterfly Arrow"

# -----------------------------------------------------------------
# Sheet "enemies": update the spawn-weight table for existing rows
# 2-8 and append new rows 9-12.
# -----------------------------------------------------------------
$ws3.Range("A2:G2").ClearContents()
$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = 0
$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 2
$ws3.Range("E2").Value = 1
$ws3.Range("F2").Value = "EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee"

$ws3.Range("A3:G3").ClearContents()
$ws3.Range("A3").Value = 2
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 0
$ws3.Range("D3").Value = 2
$ws3.Range("E3").Value = 1
$ws3.Range("F3").Value = "EnemyPrefabs/Special Enemies/Steampunk Fly/Steampunk Fly"

$ws3.Range("A4:G4").ClearContents()
$ws3.Range("A4").Value = 3
$ws3.Range("B4").Value = 0
$ws3.Range("C4").Value = 0
$ws3.Range("D4").Value = 2
$ws3.Range("E4").Value = 1
$ws3.Range("F4").Value = "EnemyPrefabs/Special Enemies/Bionic Lady Bird/Bionic Lady Bird"

$ws3.Range("A5:G5").ClearContents()
$ws3.Range("A5").Value = 4
$ws3.Range("B5").Value = 0
$ws3.Range("C5").Value = 0
$ws3.Range("D5").Value = 3
$ws3.Range("E5").Value = 1.5
$ws3.Range("F5").Value = "EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee"

$ws3.Range("A6:G6").ClearContents()
$ws3.Range("A6").Value = 5
$ws3.Range("B6").Value = 0
$ws3.Range("C6").Value = 0
$ws3.Range("D6").Value = 3
$ws3.Range("E6").Value = 1.5
$ws3.Range("F6").Value = "EnemyPrefabs/Special Enemies/Steampunk Fly/Steampunk Fly"
$ws3.Range("G6").Value = "EnemyPrefabs/Bullet Enemies//"

$ws3.Range("A7:G7").ClearContents()
$ws3.Range("A7").Value = 6
$ws3.Range("B7").Value = 0
$ws3.Range("C7").Value = 0
$ws3.Range("D7").Value = 3
$ws3.Range("E7").Value = 1.5
$ws3.Range("F7").Value = "EnemyPrefabs/Special Enemies/Bionic Lady Bird/Bionic Lady Bird"
$ws3.Range("G7").Value = "EnemyPrefabs/Special Enemies//"

$ws3.Range("A8:G8").ClearContents()
$ws3.Range("A8").Value = 7
$ws3.Range("B8").Value = 4
$ws3.Range("C8").Value = 0
$ws3.Range("D8").Value = 0
$ws3.Range("E8").Value = 1
$ws3.Range("F8").Value = "EnemyPrefabs/Arrow Enemies/Bee/Bee Arrow"

$ws3.Range("A9").Value = 8
$ws3.Range("B9").Value = 0
$ws3.Range("C9").Value = 2
$ws3.Range("D9").Value = 0
$ws3.Range("E9").Value = 1
$ws3.Range("F9").Value = "EnemyPrefabs/Bullet Enemies/Neo Fly/Neo Fly"

$ws3.Range("A10").Value = 9
$ws3.Range("B10").Value = 0
$ws3.Range("C10").Value = 0
$ws3.Range("D10").Value = 2
$ws3.Range("E10").Value = 1
$ws3.Range("F10").Value = "EnemyPrefabs/Special Enemies/Steampunk Fly/Steampunk Fly"

$ws3.Range("A11").Value = 10
$ws3.Range("B11").Value = 0
$ws3.Range("C11").Value = 0
$ws3.Range("D11").Value = 1
$ws3.Range("E11").Value = 2.5
$ws3.Range("F11").Value = "EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee"

$ws3.Range("A12").Value = 11
$ws3.Range("B12").Value = 2
$ws3.Range("C12").Value = 0
$ws3.Range("D12").Value = 0
$ws3.Range("E12").Value = 1
$ws3.Range("F12").Value = "EnemyPrefabs/Arrow Enemies/Fly/Fly Arrow"

# New rows 9-12 use the same centred style as the rest of column A.
$ws3.Range("A9:A12").HorizontalAlignment = -4108
$ws3.Range("A9:A12").VerticalAlignment = -4108

# -----------------------------------------------------------------
# Sheet "level": rewrite the per-row track counts (columns A-D) for
# the rows that changed. Column E holds formulas and is untouched.
# -----------------------------------------------------------------
$ws2.Range("A1:D1").ClearContents()
$ws2.Range("A1").Value = 11
$ws2.Range("B1").Value = 11
$ws2.Range("C1").Value = 11
$ws2.Range("D1").Value = 11

$ws2.Range("A2:D2").ClearContents()

$ws2.Range("A3:D3").ClearContents()
$ws2.Range("A3").Value = 10
$ws2.Range("B3").Value = 10
$ws2.Range("C3").Value = 10
$ws2.Range("D3").Value = 10

$ws2.Range("A5:D5").ClearContents()
$ws2.Range("A5").Value = 9
$ws2.Range("C5").Value = 7

$ws2.Range("A6:D6").ClearContents()
$ws2.Range("B6").Value = 8
$ws2.Range("D6").Value = 8

$ws2.Range("A8:D8").ClearContents()
$ws2.Range("C8").Value = 9

$ws2.Range("A9:D9").ClearContents()

$ws2.Range("A10:D10").ClearContents()
$ws2.Range("A10").Value = 11
$ws2.Range("B10").Value = 11
$ws2.Range("C10").Value = 11
$ws2.Range("D10").Value = 11

$ws2.Range("A12:D12").ClearContents()
$ws2.Range("A12").Value = 10
$ws2.Range("B12").Value = 10
$ws2.Range("C12").Value = 10
$ws2.Range("D12").Value = 10

$ws2.Range("A13:D13").ClearContents()

$ws2.Range("A14:D14").ClearContents()

$ws2.Range("A16:D16").ClearContents()
$ws2.Range("A16").Value = 8
$ws2.Range("B16").Value = 8
$ws2.Range("C16").Value = 8
$ws2.Range("D16").Value = 8

$ws2.Range("A17:D17").ClearContents()
$ws2.Range("A17").Value = 10
$ws2.Range("B17").Value = 10
$ws2.Range("C17").Value = 10
$ws2.Range("D17").Value = 10

$ws2.Range("A20:D20").ClearContents()

$ws2.Range("A23:D23").ClearContents()

Write-Host "Level data updated up to 50."
